# break out stock.yaml completed
# Convert E42:E44 (bsecode) from text to numeric values, and append three
# new rows (45-47) of breakout-screener data, mirroring rows 42-44 but for
# the 24/06/2024 11:35:50 run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")
if (-not $ws) { $ws = $wb.ActiveSheet }

# --- Re-type the existing bsecode cells (E42:E44) as numbers ---
$ws.Range("E42").Value = 20
$ws.Range("E43").Value = 531344
$ws.Range("E44").Value = 505537

# --- Append new rows 45:47 ---

# Row 45 - BSE
$ws.Range("A45").Value = "24/06/2024 11:35:50"
$ws.Range("B45").Value = 1
$ws.Range("C45").Value = "BSE"
$ws.Range("D45").Value = "BSE (Bombay stock exchange)"
$ws.Range("E45").Value = "'20"
$ws.Range("F45").Value = -2.43
$ws.Range("G45").Value = 2497.95
$ws.Range("H45").Value = 963428

# Row 46 - CONCOR
$ws.Range("A46").Value = "24/06/2024 11:35:50"
$ws.Range("B46").Value = 2
$ws.Range("C46").Value = "CONCOR"
$ws.Range("D46").Value = "Container Corporation Of India Limited"
$ws.Range("E46").Value = "'531344"
$ws.Range("F46").Value = -3.81
$ws.Range("G46").Value = 1049.35
$ws.Range("H46").Value = 5176051

# Row 47 - ZEEL
$ws.Range("A47").Value = "24/06/2024 11:35:50"
$ws.Range("B47").Value = 3
$ws.Range("C47").Value = "ZEEL"
$ws.Range("D47").Value = "Zee Entertainment Enterprises Limited"
$ws.Range("E47").Value = "'505537"
$ws.Range("F47").Value = -2.02
$ws.Range("G47").Value = 151.13
$ws.Range("H47").Value = 11607715
